$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 21 changes
$ws.Range("G2").Value = 2.4
$ws.Range("H2").Value = 3.4
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 3.4
$ws.Range("Q2").Value = 2.03
$ws.Range("R2").Value = 1.83
$ws.Range("U2").Value = 1.8
$ws.Range("V2").Value = 1.91
$ws.Range("W2").Value = 8
$ws.Range("AA2").Value = 21
$ws.Range("AB2").Value = 29
$ws.Range("AC2").Value = 10
$ws.Range("AE2").Value = 15
$ws.Range("AF2").Value = 51
$ws.Range("AG2").Value = 251
$ws.Range("AH2").Value = 9
$ws.Range("AI2").Value = 13
$ws.Range("AM2").Value = 34
$ws.Range("BD2").Value = 126

# Row 3: 8 changes
$ws.Range("G3").Value = 1.73
$ws.Range("H3").Value = 3.75
$ws.Range("L3").Value = 4.75
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("AC3").Value = 11
$ws.Range("AD3").Value = 7
$ws.Range("AG3").Value = 201

# Row 4: 18 changes
$ws.Range("G4").Value = 1.8
$ws.Range("I4").Value = 4.5
$ws.Range("J4").Value = 2.4
$ws.Range("K4").Value = 2.2
$ws.Range("L4").Value = 4.5
$ws.Range("Q4").Value = 1.9
$ws.Range("R4").Value = 1.95
$ws.Range("W4").Value = 7.5
$ws.Range("X4").Value = 8.5
$ws.Range("AB4").Value = 26
$ws.Range("AE4").Value = 15
$ws.Range("AF4").Value = 51
$ws.Range("AI4").Value = 23
$ws.Range("AN4").Value = 3.75
$ws.Range("AO4").Value = 9.5
$ws.Range("AX4").Value = 23
$ws.Range("AZ4").Value = 81
$ws.Range("BA4").Value = 101

# Row 5: 18 changes
$ws.Range("H5").Value = 3.5
$ws.Range("K5").Value = 2.1
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.4
$ws.Range("Q5").Value = 2.05
$ws.Range("R5").Value = 1.8
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 1.83
$ws.Range("W5").Value = 7
$ws.Range("X5").Value = 9.5
$ws.Range("AB5").Value = 29
$ws.Range("AC5").Value = 10
$ws.Range("AG5").Value = 301
$ws.Range("AX5").Value = 21
$ws.Range("AY5").Value = 29
$ws.Range("BC5").Value = 126

# Row 9: 28 changes
$ws.Range("J9").Value = 3.65
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 2.92
$ws.Range("Q9").Value = 1.98
$ws.Range("R9").Value = 1.75
$ws.Range("U9").Value = 1.65
$ws.Range("V9").Value = 1.98
$ws.Range("W9").Value = 9
$ws.Range("X9").Value = 16
$ws.Range("Y9").Value = 10.75
$ws.Range("AA9").Value = 28
$ws.Range("AB9").Value = 35
$ws.Range("AF9").Value = 60
$ws.Range("AH9").Value = 7.9
$ws.Range("AI9").Value = 12
$ws.Range("AL9").Value = 19
$ws.Range("AM9").Value = 27
$ws.Range("AN9").Value = 5
$ws.Range("AO9").Value = 17
$ws.Range("AP9").Value = 23
$ws.Range("AQ9").Value = 90
$ws.Range("AR9").Value = 120
$ws.Range("AS9").Value = 300
$ws.Range("AT9").Value = 2.45
$ws.Range("AX9").Value = 12.5
$ws.Range("AY9").Value = 19
$ws.Range("AZ9").Value = 50
$ws.Range("BA9").Value = 80

# Row 12: 28 changes
$ws.Range("G12").Value = 1.36
$ws.Range("H12").Value = 3.8
$ws.Range("I12").Value = 9
$ws.Range("L12").Value = 8.5
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 10
$ws.Range("Q12").Value = 2.08
$ws.Range("R12").Value = 1.73
$ws.Range("S12").Value = 1.44
$ws.Range("T12").Value = 2.63
$ws.Range("U12").Value = 2.5
$ws.Range("V12").Value = 1.5
$ws.Range("X12").Value = 5.5
$ws.Range("AA12").Value = 15
$ws.Range("AB12").Value = 41
$ws.Range("AC12").Value = 8
$ws.Range("AD12").Value = 8.5
$ws.Range("AE12").Value = 26
$ws.Range("AF12").Value = 101
$ws.Range("AH12").Value = 17
$ws.Range("AK12").Value = 126
$ws.Range("AL12").Value = 81
$ws.Range("AN12").Value = 3.1
$ws.Range("AP12").Value = 23
$ws.Range("AT12").Value = 2.63
$ws.Range("AU12").Value = 11
$ws.Range("AZ12").Value = 251
$ws.Range("BA12").Value = 301

# Row 20: 19 changes
$ws.Range("G20").Value = 1.65
$ws.Range("H20").Value = 4.2
$ws.Range("I20").Value = 5
$ws.Range("L20").Value = 5
$ws.Range("Q20").Value = 1.65
$ws.Range("R20").Value = 2.2
$ws.Range("X20").Value = 8.5
$ws.Range("AA20").Value = 12
$ws.Range("AC20").Value = 13
$ws.Range("AD20").Value = 7.5
$ws.Range("AH20").Value = 17
$ws.Range("AI20").Value = 29
$ws.Range("AL20").Value = 41
$ws.Range("AM20").Value = 41
$ws.Range("AO20").Value = 8
$ws.Range("AQ20").Value = 23
$ws.Range("AU20").Value = 8
$ws.Range("AW20").Value = 7
$ws.Range("AX20").Value = 26

# Row 21: 4 changes
$ws.Range("M21").Value = 1.03
$ws.Range("N21").Value = 17
$ws.Range("Q21").Value = 1.53
$ws.Range("R21").Value = 2.4

# Row 22: 1 changes
$ws.Range("N22").Value = 19

# Row 24: 13 changes
$ws.Range("H24").Value = 4.1
$ws.Range("I24").Value = 3.3
$ws.Range("J24").Value = 2.4
$ws.Range("K24").Value = 2.6
$ws.Range("AC24").Value = 23
$ws.Range("AH24").Value = 19
$ws.Range("AI24").Value = 23
$ws.Range("AM24").Value = 21
$ws.Range("AS24").Value = 67
$ws.Range("AU24").Value = 6.5
$ws.Range("AV24").Value = 34
$ws.Range("BC24").Value = 201
$ws.Range("BD24").Value = 151

# Row 26: 19 changes
$ws.Range("G26").Value = 2.2
$ws.Range("I26").Value = 3.1
$ws.Range("J26").Value = 2.75
$ws.Range("K26").Value = 2.25
$ws.Range("L26").Value = 3.6
$ws.Range("U26").Value = 1.67
$ws.Range("V26").Value = 2.1
$ws.Range("W26").Value = 9
$ws.Range("X26").Value = 11
$ws.Range("AK26").Value = 34
$ws.Range("AL26").Value = 23
$ws.Range("AM26").Value = 29
$ws.Range("AN26").Value = 4.33
$ws.Range("AP26").Value = 21
$ws.Range("AV26").Value = 51
$ws.Range("AX26").Value = 17
$ws.Range("AY26").Value = 23
$ws.Range("BB26").Value = 151
$ws.Range("BC26").Value = 501

# Row 28: 1 changes
$ws.Range("N28").Value = 19

# Row 32: 13 changes
$ws.Range("G32").Value = 2.3
$ws.Range("I32").Value = 3
$ws.Range("K32").Value = 2.2
$ws.Range("L32").Value = 3.5
$ws.Range("O32").Value = 1.25
$ws.Range("P32").Value = 3.75
$ws.Range("Q32").Value = 1.83
$ws.Range("R32").Value = 2.03
$ws.Range("U32").Value = 1.67
$ws.Range("V32").Value = 2.1
$ws.Range("W32").Value = 9
$ws.Range("Z32").Value = 21
$ws.Range("AL32").Value = 23
